$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update underlying grade entries (formulas in row 7 / J7 recalc automatically)
$ws.Range("G9").Value = 0.675
$ws.Range("G10").Value = 0.6
$ws.Range("F13").Value = 0.5

# Move the active selection to match the author's final cursor position
$ws.Activate()
$ws.Range("G10").Select() | Out-Null
